$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2060810810810811
$ws.Range("C2").Value = 0.5337837837837838
$ws.Range("J2").Value = 0.01013513513513514
$ws.Range("P2").Value = 0.1621621621621622
$ws.Range("S2").Value = 0.08783783783783784

# Row 3
$ws.Range("B3").Value = 0.01219512195121951
$ws.Range("C3").Value = 0.03658536585365853
$ws.Range("J3").Value = 0.06097560975609756
$ws.Range("P3").Value = 0.75
$ws.Range("S3").Value = 0.1402439024390244

# Row 4
$ws.Range("J4").Value = 0.03448275862068965
$ws.Range("O4").Value = 0.03448275862068965
$ws.Range("P4").Value = 0.6206896551724138
$ws.Range("S4").Value = 0.3103448275862069

# Row 5
$ws.Range("P5").Value = 0.5
$ws.Range("S5").Value = 0.5

# Row 6
$ws.Range("B6").Value = 0.08947368421052632
$ws.Range("D6").Value = 0.005263157894736842
$ws.Range("E6").Value = 0.005263157894736842
$ws.Range("F6").Value = 0.06842105263157895
$ws.Range("J6").Value = 0.3263157894736842
$ws.Range("O6").Value = 0.02631578947368421
$ws.Range("Q6").Value = 0.07368421052631578
$ws.Range("R6").Value = 0.07368421052631578
$ws.Range("S6").Value = 0.331578947368421

# Row 7
$ws.Range("B7").Value = 0.1344086021505376
$ws.Range("D7").Value = 0.01612903225806452
$ws.Range("E7").Value = 0.005376344086021506
$ws.Range("F7").Value = 0.06989247311827956
$ws.Range("J7").Value = 0.1612903225806452
$ws.Range("O7").Value = 0.01612903225806452
$ws.Range("Q7").Value = 0.1236559139784946
$ws.Range("R7").Value = 0.06989247311827956
$ws.Range("S7").Value = 0.4032258064516129

# Row 8
$ws.Range("B8").Value = 0.1253196930946291
$ws.Range("D8").Value = 0.01534526854219949
$ws.Range("F8").Value = 0.07416879795396419
$ws.Range("J8").Value = 0.09974424552429667
$ws.Range("O8").Value = 0.01790281329923274
$ws.Range("Q8").Value = 0.1253196930946291
$ws.Range("R8").Value = 0.1150895140664962
$ws.Range("S8").Value = 0.4271099744245525

# Row 9
$ws.Range("B9").Value = 0.08943089430894309
$ws.Range("D9").Value = 0.01626016260162602
$ws.Range("F9").Value = 0.06504065040650407
$ws.Range("J9").Value = 0.08130081300813008
$ws.Range("O9").Value = 0.02439024390243903
$ws.Range("Q9").Value = 0.1544715447154472
$ws.Range("R9").Value = 0.08943089430894309
$ws.Range("S9").Value = 0.4796747967479675

# Row 10
$ws.Range("B10").Value = 0.1289340101522843
$ws.Range("D10").Value = 0.02233502538071066
$ws.Range("E10").Value = 0.002030456852791878
$ws.Range("F10").Value = 0.05583756345177665
$ws.Range("J10").Value = 0.1157360406091371
$ws.Range("O10").Value = 0.01218274111675127
$ws.Range("Q10").Value = 0.1644670050761421
$ws.Range("R10").Value = 0.0751269035532995
$ws.Range("S10").Value = 0.4233502538071066

# Row 11
$ws.Range("G11").Value = 0.1164383561643836
$ws.Range("J11").Value = 0.08561643835616438
$ws.Range("K11").Value = 0.1815068493150685
$ws.Range("L11").Value = 0.5958904109589042
$ws.Range("S11").Value = 0.02054794520547945

# Row 12
$ws.Range("G12").Value = 0.7428571428571429
$ws.Range("J12").Value = 0.1885714285714286
$ws.Range("K12").Value = 0.005714285714285714
$ws.Range("L12").Value = 0.01142857142857143
$ws.Range("S12").Value = 0.05142857142857143

# Row 13
$ws.Range("G13").Value = 0.6857142857142857
$ws.Range("J13").Value = 0.2857142857142857
$ws.Range("S13").Value = 0.02857142857142857

# Row 14
$ws.Range("G14").Value = 0.5
$ws.Range("S14").Value = 0.5

# Row 15
$ws.Range("F15").Value = 0.02469135802469136
$ws.Range("H15").Value = 0.1604938271604938
$ws.Range("I15").Value = 0.04938271604938271
$ws.Range("J15").Value = 0.4197530864197531
$ws.Range("K15").Value = 0.04938271604938271
$ws.Range("M15").Value = 0.006172839506172839
$ws.Range("O15").Value = 0.04938271604938271
$ws.Range("S15").Value = 0.2407407407407407

# Row 16
$ws.Range("F16").Value = 0.0481283422459893
$ws.Range("H16").Value = 0.106951871657754
$ws.Range("I16").Value = 0.05882352941176471
$ws.Range("J16").Value = 0.4171122994652406
$ws.Range("K16").Value = 0.160427807486631
$ws.Range("M16").Value = 0.0160427807486631
$ws.Range("O16").Value = 0.053475935828877
$ws.Range("S16").Value = 0.1390374331550802

# Row 17
$ws.Range("F17").Value = 0.04150943396226415
$ws.Range("H17").Value = 0.1584905660377358
$ws.Range("I17").Value = 0.0830188679245283
$ws.Range("J17").Value = 0.369811320754717
$ws.Range("K17").Value = 0.1509433962264151
$ws.Range("M17").Value = 0.01132075471698113
$ws.Range("O17").Value = 0.05660377358490566
$ws.Range("S17").Value = 0.1283018867924528

# Row 18
$ws.Range("F18").Value = 0.03125
$ws.Range("H18").Value = 0.24375
$ws.Range("I18").Value = 0.1
$ws.Range("J18").Value = 0.35
$ws.Range("K18").Value = 0.06875000000000001
$ws.Range("M18").Value = 0.03125
$ws.Range("N18").Value = 0.00625
$ws.Range("O18").Value = 0.0625
$ws.Range("S18").Value = 0.10625

# Row 19
$ws.Range("F19").Value = 0.01886792452830189
$ws.Range("H19").Value = 0.2345013477088949
$ws.Range("I19").Value = 0.06199460916442048
$ws.Range("J19").Value = 0.3324348607367476
$ws.Range("K19").Value = 0.12848158131177
$ws.Range("M19").Value = 0.0215633423180593
$ws.Range("N19").Value = 0.001796945193171608
$ws.Range("O19").Value = 0.06019766397124888
$ws.Range("S19").Value = 0.1401617250673854

Write-Output "Applied updated simulation probabilities"